$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16:D16").Value = "Create Email"
$ws.Range("A17:D17").Value = "Email Members"

$ws.Range("A1").Select()
